# Apply cached-value corrections to the profit-calculation columns (H:N)
# across all item sheets, per the scheduled pricing-data refresh.
$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 1333.3334
$ws.Range("J18").Value = 1000
$ws.Range("L18").Value = 1000
$ws.Range("N18").Value = -1568
$ws.Range("H19").Value = 562.94446
$ws.Range("I19").Value = 520.1111
$ws.Range("K19").Value = 520.1111
$ws.Range("M19").Value = -345.1111
$ws.Range("H132").Value = 2015.1724
$ws.Range("I132").Value = 1692.5217
$ws.Range("K132").Value = 5077.5651
$ws.Range("M132").Value = -2547.5651
$ws.Range("H138").Value = 3056.3115
$ws.Range("I138").Value = 1639.1538
$ws.Range("J138").Value = 4109.057
$ws.Range("K138").Value = 4917.4614
$ws.Range("L138").Value = 12327.171
$ws.Range("M138").Value = 222.5385999999999
$ws.Range("N138").Value = -22607.171
# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 3318.3333
$ws.Range("I5").Value = 3318.3333
$ws.Range("K5").Value = 3318.3333
$ws.Range("M5").Value = -3206.3333
$ws.Range("H32").Value = 13683.482
$ws.Range("I32").Value = 15407.149
$ws.Range("K32").Value = 15407.149
$ws.Range("M32").Value = -15120.149
$ws.Range("H61").Value = 2299.5417
$ws.Range("I61").Value = 1799.421
$ws.Range("J61").Value = 4200
$ws.Range("K61").Value = 1799.421
$ws.Range("L61").Value = 4200
$ws.Range("M61").Value = -1587.421
$ws.Range("N61").Value = -4624
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()
$ws.Range("H136").Value = 2299.5417
$ws.Range("I136").Value = 1799.421
$ws.Range("J136").Value = 4200
$ws.Range("K136").Value = 5398.263
$ws.Range("L136").Value = 12600
$ws.Range("M136").Value = -2848.263
$ws.Range("N136").Value = -17700
# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 3318.3333
$ws.Range("I4").Value = 3318.3333
$ws.Range("K4").Value = 3318.3333
$ws.Range("M4").Value = -3203.3333
$ws.Range("H134").Value = 2759.7334
$ws.Range("I134").Value = 2371.524
$ws.Range("J134").Value = 3665.5557
$ws.Range("K134").Value = 7114.572
$ws.Range("L134").Value = 10996.6671
$ws.Range("M134").Value = -4579.572
$ws.Range("N134").Value = -16066.6671
# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 5014726
$ws.Range("I6").Value = 10000001
$ws.Range("K6").Value = 10000001
$ws.Range("M6").Value = -9999888
$ws.Range("H7").Value = 130.83333
$ws.Range("I7").Value = 100
$ws.Range("J7").Value = 137
$ws.Range("K7").Value = 100
$ws.Range("L7").Value = 137
$ws.Range("M7").Value = 13
$ws.Range("N7").Value = -363
$ws.Range("H41").Value = 5000
$ws.Range("I41").Value = 5000
$ws.Range("K41").Value = 5000
$ws.Range("M41").Value = -4572
$ws.Range("H50").Value = 24850
$ws.Range("J50").Value = 24850
$ws.Range("L50").Value = 24850
$ws.Range("N50").Value = -26100
$ws.Range("H51").Value = 18940.2
$ws.Range("I51").Value = 16000.333
$ws.Range("J51").Value = 23350
$ws.Range("K51").Value = 16000.333
$ws.Range("L51").Value = 23350
$ws.Range("M51").Value = -15264.333
$ws.Range("N51").Value = -24822
$ws.Range("H59").Value = 27806.25
$ws.Range("J59").Value = 27806.25
$ws.Range("L59").Value = 27806.25
$ws.Range("N59").Value = -30096.25
$ws.Range("H61").Value = 18940.2
$ws.Range("I61").Value = 16000.333
$ws.Range("J61").Value = 23350
$ws.Range("K61").Value = 16000.333
$ws.Range("L61").Value = 23350
$ws.Range("M61").Value = -15652.333
$ws.Range("N61").Value = -24046
$ws.Range("H74").Value = 34104.668
$ws.Range("J74").Value = 34104.668
$ws.Range("L74").Value = 34104.668
$ws.Range("N74").Value = -35852.668
$ws.Range("H77").Value = 34104.668
$ws.Range("J77").Value = 34104.668
$ws.Range("L77").Value = 102314.004
$ws.Range("N77").Value = -111050.004
$ws.Range("H132").Value = 2531
$ws.Range("I132").Value = 1480.2084
$ws.Range("J132").Value = 5333.1113
$ws.Range("K132").Value = 4440.6252
$ws.Range("L132").Value = 15999.3339
$ws.Range("M132").Value = -1910.6252
$ws.Range("N132").Value = -21059.3339
# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 6030.3
$ws.Range("I80").Value = 11500
$ws.Range("J80").Value = 3686.1428
$ws.Range("K80").Value = 34500
$ws.Range("L80").Value = 11058.4284
$ws.Range("M80").Value = -33564
$ws.Range("N80").Value = -12930.4284
$ws.Range("H83").Value = 6030.3
$ws.Range("I83").Value = 11500
$ws.Range("J83").Value = 3686.1428
$ws.Range("K83").Value = 103500
$ws.Range("L83").Value = 33175.2852
$ws.Range("M83").Value = -98820
$ws.Range("N83").Value = -42535.2852
# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 841837.5
$ws.Range("I2").Value = 1262675
$ws.Range("J2").Value = 162.5
$ws.Range("K2").Value = 1262675
$ws.Range("L2").Value = 162.5
$ws.Range("M2").Value = -1262562
$ws.Range("N2").Value = -388.5
$ws.Range("H102").Value = 2849.375
$ws.Range("I102").Value = 2812.4443
$ws.Range("K102").Value = 2812.4443
$ws.Range("M102").Value = -1190.4443
$ws.Range("H132").Value = 3830.7
$ws.Range("I132").Value = 3042.3
$ws.Range("J132").Value = 4224.9
$ws.Range("K132").Value = 9126.900000000001
$ws.Range("L132").Value = 12674.7
$ws.Range("M132").Value = -6596.900000000001
$ws.Range("N132").Value = -17734.7
# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 9401.143
$ws.Range("I2").Value = 9000
$ws.Range("K2").Value = 9000
$ws.Range("M2").Value = -8888
$ws.Range("H122").Value = 32146714
$ws.Range("I122").Value = 25003938
$ws.Range("J122").Value = 50003652
$ws.Range("K122").Value = 75011814
$ws.Range("L122").Value = 150010956
$ws.Range("M122").Value = -75009364
$ws.Range("N122").Value = -150015856
# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 40000
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 40000
$ws.Range("K54").Value = 0
$ws.Range("L54").Value = 40000
$ws.Range("M54").ClearContents()
$ws.Range("N54").Value = -41040
$ws.Range("H132").Value = 3287.4
$ws.Range("I132").Value = 3627.25
$ws.Range("J132").Value = 3163.818
$ws.Range("K132").Value = 10881.75
$ws.Range("L132").Value = 9491.454000000002
$ws.Range("M132").Value = -8351.75
$ws.Range("N132").Value = -14551.454
